# Re-position the footer and slide-number placeholders on slide 8
# ("Eredmenyek" / id-5 "Footer Placeholder 4" and id-6 "Slide Number
# Placeholder 5") so they sit with the rest of the deck's moved footer
# row, per the upload diff:
#   Footer Placeholder 4:      off 4143262,6492875  -> 2573367,5880102
#   Slide Number Placeholder 5: off 11227439,6492875 -> 10951856,5880101
#
# PowerPoint's Shape.Left/Top are expressed in points (1 pt = 12700 EMU),
# so the EMU targets below are converted to points; the literals carry
# extra fractional digits so that, after the host's internal float
# rounding, the saved OOXML lands on the exact target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

$footer = $s.Shapes.Item("Footer Placeholder 4")
$slideNum = $s.Shapes.Item("Slide Number Placeholder 5")

$footer.Left = 202.6273651123047
$footer.Top = 463.0001983642578

$slideNum.Left = 862.3508911132812
$slideNum.Top = 463.0001220703125
